$wb = $excel.ActiveWorkbook

# --- BBNPPTY sheet: change the "111 rules for new gas" representation ---
# Instead of banning new combined-cycle gas capacity from 2032 onward (values
# of 1 in M4:AE4), represent it as no change in new capacity factor (values
# of 0 in M4:AE4).
$wsData = $wb.Worksheets.Item("BBNPPTY")
$wsData.Range("M4:AE4").Value = 0

# --- About sheet: update the notes text ---
# The note cell A13 used to hold descriptive text ("2028 and new combined
# cycle gas without CCS is banned starting in 2032.") which is no longer
# accurate now that the ban has been replaced by a capacity-factor change.
# Replace it with the literal year value that is still relevant (2028).
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A13").Value = 2028

# --- Restore selections/active cells to match the saved view state ---
# BBNPPTY: selection over the 2031 (L) through 2050 (AE) columns of row 4.
$wsData.Select()
$wsData.Range("L4:AE4").Select()

# About: leave the sheet selected on cell B18, and make sure "About" ends
# up as the active sheet/tab (it was the tab selected before the edit).
$wsAbout.Select()
$wsAbout.Range("B18").Select()
